$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 = "CHMO_DEF"
$ws.Range('F1').Value = 'CHMO_DEF'

# Fill column F (CHMO_DEF) for data rows 2-13 with definition text
$ws.Range('F2').Value = '[''A material entity is an independent continuant that has some portion of matter as proper or improper continuant part. [BFO]'']'
$ws.Range('F3').Value = '[''p is a process if p is an occurrent that has temporal proper parts and for some time t, p specifically depends on some material entity at t. [BFO]'', locstr("Process, i.e., a physical entity with a temporal evolution that ''has a meaning for the ontologist''", ''en'')]'
$ws.Range('F4').Value = '[''A function is a disposition that exists in virtue of the bearer’s physical make-up and this physical make-up is something the bearer possesses because it came into being, either through evolution (in the case of natural biological entities) or through intentional design (in the case of artifacts), in order to realize processes of a certain sort. [BFO]'']'
$ws.Range('F5').Value = '[''A quality is a specifically dependent continuant that, in contrast to roles and dispositions, does not require any further process in order to be realized. [BFO]'']'
$ws.Range('F6').Value = '[''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]'', ''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])'']'
$ws.Range('F7').Value = '[''A continuant is an entity that persists, endures, or continues to exist through time while maintaining its identity. [BFO]'']'
$ws.Range('F8').Value = '[''b is an independent continuant if b is a continuant which is such that there is no c and no t such that b s-depends_on c at t. [BFO]'']'
$ws.Range('F9').Value = '[''b is a specifically dependent continuant if b is a continuant and there is some independent continuant c which is not a spatial region and which is such that b specifically depends on c at every time t during the course of b’s existence. [BFO]'']'
$ws.Range('F10').Value = '[''b is a generically dependent continuant if b is a continuant that generically depends on one or more other entities. [BFO]'']'
$ws.Range('F11').Value = '[''B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]'']'
$ws.Range('F12').Value = '[''B is a disposition means: b is a realizable entity and b’s bearer is some material entity and b is such that if it ceases to exist, then its bearer is physically changed, and b’s realization occurs when and because this bearer is in some special physical circumstances, and this realization occurs in virtue of the bearer’s physical make-up. [BFO]'']'
$ws.Range('F13').Value = '[''B is a continuant fiat boundary if b is an immaterial entity that is of zero, one or two dimensions and does not include a spatial region as part. [BFO]'']'

Write-Output "Done updating CHMO_DEF column"
